$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# erros (column G) increments
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1

# acertos (column F) increment
$ws.Range("F6").Value = 2

# erros (column G) increment
$ws.Range("G7").Value = 2
